# Deploy updated output folder
# Update the "Title" and "Date" metadata values on the Metadata sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 5: Title -> "NG-Imm Sibling Health Status VS"
$ws.Range("B5").Value = "NG-Imm Sibling Health Status VS"

# Row 8: Date -> "2025-06-24T09:13:37+01:00"
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
